$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the formatting already used by sibling cells in the "Location" column
# (style index 1 / General) before writing the geocoding text, so the F cells
# stop looking like the inherited (empty) date-formatted cells.
$ws.Range("C2").Copy()
$ws.Range("F2:F6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F8:F10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("F2").Value = "MacDonalds Restaurant, Craigieburn"
$ws.Range("F3").Value = "Myer Highpoint, Maribyrnong"
$ws.Range("F4").Value = "MacDonalds Restaurant, Fawkner"
$ws.Range("F5").Value = "Grant Lodge, Bacchus Marsh"
$ws.Range("F6").Value = "Cedar Meats Australia, Brooklyn "
$ws.Range("F8").Value = "The Learning Sanctuary, Yarraville"
$ws.Range("F9").Value = "Doutta Galla Aged Care Home, Footscray"
$ws.Range("F10").Value = "Sunshine Hospital, Sunshine"

# Reflect the selection/scroll state captured in the saved workbook: the
# user scrolled the sheet so column C is the first visible column, and had
# F2:F6 selected with F2 as the active cell.
$ws.Range("F2:F6").Select()
$ws.Range("F2").Activate()
$excel.ActiveWindow.ScrollColumn = 3
